# Update gh-pages to output generated at 456a3b4
# This script updates the "想去人数" (want-to-go count) column F values
# across the 展览, 演出 and 全部类型 worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 509
$ws1.Range("F5").Value = 147
$ws1.Range("F7").Value = 298
$ws1.Range("F8").Value = 5210
$ws1.Range("F9").Value = 137
$ws1.Range("F10").Value = 725
$ws1.Range("F13").Value = 317
$ws1.Range("F14").Value = 47
$ws1.Range("F15").Value = 6332
$ws1.Range("F18").Value = 143
$ws1.Range("F20").Value = 15177
$ws1.Range("F21").Value = 1506
$ws1.Range("F22").Value = 272
$ws1.Range("F23").Value = 132
$ws1.Range("F25").Value = 10998
$ws1.Range("F26").Value = 737
$ws1.Range("F27").Value = 4289
$ws1.Range("F28").Value = 224
$ws1.Range("F29").Value = 371
$ws1.Range("F30").Value = 10

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 42

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 509
$ws4.Range("F5").Value = 147
$ws4.Range("F7").Value = 298
$ws4.Range("F9").Value = 5210
$ws4.Range("F10").Value = 137
$ws4.Range("F11").Value = 725
$ws4.Range("F15").Value = 317
$ws4.Range("F16").Value = 47
$ws4.Range("F17").Value = 42
$ws4.Range("F18").Value = 6333
$ws4.Range("F21").Value = 143
$ws4.Range("F23").Value = 15177
$ws4.Range("F24").Value = 1506
$ws4.Range("F25").Value = 272
$ws4.Range("F26").Value = 132
$ws4.Range("F28").Value = 10998
$ws4.Range("F29").Value = 737
$ws4.Range("F30").Value = 4289
$ws4.Range("F31").Value = 224
$ws4.Range("F32").Value = 371
$ws4.Range("F33").Value = 10
